$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 354
$ws.Range("I2").Value = 338.83334
$ws.Range("J2").Value = 445
$ws.Range("K2").Value = 338.83334
$ws.Range("L2").Value = 445
$ws.Range("M2").Value = -225.83334
$ws.Range("N2").Value = -671

$ws.Range("H64").Value = 4750
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 4750
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 4750
$ws.Range("N64").Value = -5246

$ws.Range("H67").Value = 4750
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 4750
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 4750
$ws.Range("N67").Value = -6466

$ws.Range("H100").Value = 1652.7273
$ws.Range("I100").Value = 1576.1111
$ws.Range("J100").Value = 1997.5
$ws.Range("K100").Value = 1576.1111
$ws.Range("L100").Value = 1997.5
$ws.Range("M100").Value = -1035.1111

$ws.Range("H132").Value = 2981.4194
$ws.Range("I132").Value = 2869.6206
$ws.Range("J132").Value = 4602.5
$ws.Range("K132").Value = 8608.861800000001
$ws.Range("L132").Value = 13807.5
$ws.Range("M132").Value = -6078.861800000001

$ws.Range("H137").Value = 62501050
$ws.Range("I137").Value = 83334310
$ws.Range("J137").Value = 1250
$ws.Range("K137").Value = 250002930
$ws.Range("L137").Value = 3750
$ws.Range("M137").Value = -250000380

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1316.4062
$ws.Range("I32").Value = 1316.4062
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1316.4062
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1029.4062

$ws.Range("H60").Value = 17666.334
$ws.Range("I60").Value = 17666.334
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 17666.334
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -16933.334
$ws.Range("N60").ClearContents()

$ws.Range("H61").Value = 12821870
$ws.Range("I61").Value = 13334706
$ws.Range("J61").Value = 980
$ws.Range("K61").Value = 13334706
$ws.Range("L61").Value = 980
$ws.Range("M61").Value = -13334494
$ws.Range("N61").Value = -1404

$ws.Range("H132").Value = 71430830
$ws.Range("I132").Value = 83335496
$ws.Range("J132").Value = 2850
$ws.Range("K132").Value = 250006488
$ws.Range("L132").Value = 8550
$ws.Range("M132").Value = -250003958
$ws.Range("N132").Value = -13610

$ws.Range("H136").Value = 12821870
$ws.Range("I136").Value = 13334706
$ws.Range("J136").Value = 980
$ws.Range("K136").Value = 40004118
$ws.Range("L136").Value = 2940
$ws.Range("M136").Value = -40001568
$ws.Range("N136").Value = -8040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 22729688
$ws.Range("I86").Value = 29414232
$ws.Range("J86").Value = 2235.2
$ws.Range("K86").Value = 29414232
$ws.Range("L86").Value = 2235.2
$ws.Range("M86").Value = -29413109

$ws.Range("H89").Value = 22729688
$ws.Range("I89").Value = 29414232
$ws.Range("J89").Value = 2235.2
$ws.Range("K89").Value = 147071160
$ws.Range("L89").Value = 11176
$ws.Range("M89").Value = -147065544

$ws.Range("H99").Value = 1181.6
$ws.Range("I99").Value = 1099.3334
$ws.Range("J99").Value = 1305
$ws.Range("K99").Value = 1099.3334
$ws.Range("L99").Value = 1305
$ws.Range("M99").Value = 398.6666

$ws.Range("H107").Value = 15627597
$ws.Range("I107").Value = 2754.6072
$ws.Range("J107").Value = 125001496
$ws.Range("K107").Value = 2754.6072
$ws.Range("L107").Value = 125001496
$ws.Range("M107").Value = -834.6071999999999

$ws.Range("H134").Value = 3019.5356
$ws.Range("I134").Value = 3006.0833
$ws.Range("J134").Value = 3100.25
$ws.Range("K134").Value = 9018.249899999999
$ws.Range("L134").Value = 9300.75
$ws.Range("M134").Value = -6483.249899999999
$ws.Range("N134").Value = -14370.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1614.0834
$ws.Range("I58").Value = 1252.2858
$ws.Range("J58").Value = 2120.6
$ws.Range("K58").Value = 1252.2858
$ws.Range("L58").Value = 2120.6
$ws.Range("M58").Value = -1049.2858
$ws.Range("N58").Value = -2526.6

$ws.Range("H60").Value = 19350.5
$ws.Range("I60").Value = 17333.334
$ws.Range("J60").Value = 21367.666
$ws.Range("K60").Value = 17333.334
$ws.Range("L60").Value = 21367.666
$ws.Range("M60").Value = -16822.334
$ws.Range("N60").Value = -22389.666

$ws.Range("H94").Value = 760
$ws.Range("I94").Value = 732.4545000000001
$ws.Range("J94").Value = 783.3077
$ws.Range("K94").Value = 732.4545000000001
$ws.Range("L94").Value = 783.3077
$ws.Range("M94").Value = -281.4545000000001
$ws.Range("N94").Value = -1685.3077

$ws.Range("H132").Value = 2717.8518
$ws.Range("I132").Value = 2464.44
$ws.Range("J132").Value = 5885.5
$ws.Range("K132").Value = 7393.32
$ws.Range("L132").Value = 17656.5
$ws.Range("M132").Value = -4863.32

$ws.Range("H134").Value = 1879.6
$ws.Range("I134").Value = 1539.5333
$ws.Range("J134").Value = 2899.8
$ws.Range("K134").Value = 4618.5999
$ws.Range("L134").Value = 8699.400000000001
$ws.Range("M134").Value = -2083.5999

$ws.Range("H136").Value = 1614.0834
$ws.Range("I136").Value = 1252.2858
$ws.Range("J136").Value = 2120.6
$ws.Range("K136").Value = 3756.8574
$ws.Range("L136").Value = 6361.799999999999
$ws.Range("M136").Value = -1206.8574
$ws.Range("N136").Value = -11461.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 336.66666
$ws.Range("I117").Value = 195
$ws.Range("J117").Value = 407.5
$ws.Range("K117").Value = 585
$ws.Range("L117").Value = 1222.5
$ws.Range("M117").Value = 2857
$ws.Range("N117").Value = -8106.5

$ws.Range("H137").Value = 2140
$ws.Range("I137").Value = 900
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 2700
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = 2400

$ws.Range("H140").Value = 1671
$ws.Range("I140").Value = 1401.25
$ws.Range("J140").Value = 2750
$ws.Range("K140").Value = 4203.75
$ws.Range("L140").Value = 8250
$ws.Range("M140").Value = 976.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 81.75
$ws.Range("I2").Value = 21.166666
$ws.Range("J2").Value = 142.33333
$ws.Range("K2").Value = 21.166666
$ws.Range("L2").Value = 142.33333
$ws.Range("M2").Value = 91.83333400000001
$ws.Range("N2").Value = -368.33333

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()

$ws.Range("H132").Value = 2913.5813
$ws.Range("I132").Value = 2086.76
$ws.Range("J132").Value = 4061.9443
$ws.Range("K132").Value = 6260.280000000001
$ws.Range("L132").Value = 12185.8329
$ws.Range("M132").Value = -3730.280000000001
$ws.Range("N132").Value = -17245.8329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 166667330
$ws.Range("I7").Value = 250000500
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 250000500
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -250000388

$ws.Range("H46").Value = 1570.7778
$ws.Range("I46").Value = 863.5
$ws.Range("J46").Value = 1986.8235
$ws.Range("K46").Value = 863.5
$ws.Range("L46").Value = 1986.8235
$ws.Range("M46").Value = -675.5
$ws.Range("N46").Value = -2362.8235

$ws.Range("H55").Value = 266.14285
$ws.Range("I55").Value = 262.2143
$ws.Range("J55").Value = 274
$ws.Range("K55").Value = 262.2143
$ws.Range("L55").Value = 274
$ws.Range("M55").Value = -89.21429999999998

$ws.Range("H61").Value = 14143
$ws.Range("I61").Value = 10316.667
$ws.Range("J61").Value = 17422.715
$ws.Range("K61").Value = 10316.667
$ws.Range("L61").Value = 17422.715
$ws.Range("M61").Value = -10114.667
$ws.Range("N61").Value = -17826.715

$ws.Range("H113").Value = 14143
$ws.Range("I113").Value = 10316.667
$ws.Range("J113").Value = 17422.715
$ws.Range("K113").Value = 10316.667
$ws.Range("L113").Value = 17422.715
$ws.Range("M113").Value = -8146.666999999999
$ws.Range("N113").Value = -21762.715

$ws.Range("H126").Value = 166667330
$ws.Range("I126").Value = 250000500
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 750001500
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -749999030

$ws.Range("H132").Value = 2801.2683
$ws.Range("I132").Value = 2374.36
$ws.Range("J132").Value = 3468.3125
$ws.Range("K132").Value = 7123.08
$ws.Range("L132").Value = 10404.9375
$ws.Range("M132").Value = -4593.08

$ws.Range("H136").Value = 1947.5106
$ws.Range("I136").Value = 1850.5151
$ws.Range("J136").Value = 2176.1428
$ws.Range("K136").Value = 5551.5453
$ws.Range("L136").Value = 6528.428400000001
$ws.Range("M136").Value = -3001.5453

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2000
$ws.Range("N5").Value = -2224
$ws.Range("M5").ClearContents()

$ws.Range("H64").Value = 29998
$ws.Range("I64").Value = 29998
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 29998
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -29750
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 29998
$ws.Range("I67").Value = 29998
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 29998
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -29140
$ws.Range("N67").ClearContents()

$ws.Range("H119").Value = 45037
$ws.Range("I119").Value = 52500
$ws.Range("J119").Value = 42549.332
$ws.Range("K119").Value = 52500
$ws.Range("L119").Value = 42549.332
$ws.Range("M119").Value = -47662
$ws.Range("N119").Value = -52225.332

$ws.Range("H122").Value = 1877.0769
$ws.Range("I122").Value = 1943.1428
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5829.428400000001
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -3379.428400000001
$ws.Range("N122").Value = -10300

$ws.Range("H132").Value = 6423.3125
$ws.Range("I132").Value = 6930.4443
$ws.Range("J132").Value = 5771.2856
$ws.Range("K132").Value = 20791.3329
$ws.Range("L132").Value = 17313.8568
$ws.Range("M132").Value = -18261.3329

$ws.Range("H136").Value = 3382.2222
$ws.Range("I136").Value = 952.53845
$ws.Range("J136").Value = 9699.4
$ws.Range("K136").Value = 2857.61535
$ws.Range("L136").Value = 29098.2
$ws.Range("M136").Value = -307.61535
$ws.Range("N136").Value = -34198.2
